# Update "想去人数" (F column) counts for the 展览 (Exhibitions) and
# 全部类型 (All types) sheets, which share identical data.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 858
    "F4"  = 2172
    "F6"  = 12626
    "F7"  = 59
    "F11" = 1151
    "F12" = 954
    "F13" = 13665
    "F14" = 14003
    "F19" = 16
    "F26" = 707
    "F27" = 5122
    "F28" = 5
    "F29" = 255
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
